$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure columns B:E stay formatted as plain text so values
# (e.g. "28.487.33", "0.9992", "0.00001015") are written verbatim
# instead of being reinterpreted as numbers by Excel.
$ws.Range("B2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '28.487.33'
$ws.Range("E2").Value = '  +0.51%  '
$ws.Range("D3").Value = '1.864.92'
$ws.Range("E3").Value = '  -0.48%  '
$ws.Range("D4").Value = '0.9992'
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").Value = '328.94'
$ws.Range("E5").Value = '  -2.46%  '
$ws.Range("D6").Value = '0.9975'
$ws.Range("E6").Value = '  -0.16%  '
$ws.Range("D7").Value = '0.4569'
$ws.Range("E7").Value = '  -2.89%  '
$ws.Range("D8").Value = '0.3950'
$ws.Range("E8").Value = '  +0.47%  '
$ws.Range("D9").Value = '47.63'
$ws.Range("E9").Value = '  +1.17%  '
$ws.Range("D10").Value = '0.07834'
$ws.Range("E10").Value = '  -1.95%  '
$ws.Range("D11").Value = '0.9848'
$ws.Range("E11").Value = '  -2.99%  '
$ws.Range("D12").Value = '21.37'
$ws.Range("E12").Value = '  -1.99%  '
$ws.Range("D13").Value = '1.851.54'
$ws.Range("E13").Value = '  -0.15%  '
$ws.Range("E14").Value = '  -2.96%  '
$ws.Range("D15").Value = '6.982'
$ws.Range("E15").Value = '  -4.07%  '
$ws.Range("D16").Value = '0.9998'
$ws.Range("E16").Value = '  -0.07%  '
$ws.Range("D17").Value = '88.02'
$ws.Range("E17").Value = '  -3.60%  '
$ws.Range("D18").Value = '0.06521'
$ws.Range("E18").Value = '  -1.07%  '
$ws.Range("D19").Value = '0.00001015'
$ws.Range("E19").Value = '  -2.69%  '
$ws.Range("D20").Value = '17.11'
$ws.Range("E20").Value = '  -3.58%  '
$ws.Range("D21").Value = '1.005'
$ws.Range("E21").Value = '  +0.47%  '
$ws.Range("D22").Value = '28.482.27'
$ws.Range("E22").Value = '  +0.54%  '
$ws.Range("D23").Value = '5.312'
$ws.Range("E23").Value = '  -2.70%  '
$ws.Range("D24").Value = '10.79'
$ws.Range("E24").Value = '  -2.58%  '
$ws.Range("D25").Value = '2.246'
$ws.Range("E25").Value = '  -1.78%  '
$ws.Range("D26").Value = '2.083.28'
$ws.Range("E26").Value = '  +0.48%  '
$ws.Range("D27").Value = '156.77'
$ws.Range("E27").Value = '  -1.34%  '
$ws.Range("D28").Value = '19.23'
$ws.Range("E28").Value = '  -3.86%  '
$ws.Range("D29").Value = '2.070'
$ws.Range("E29").Value = '  -4.05%  '
$ws.Range("D30").Value = '5.302'
$ws.Range("E30").Value = '  -4.09%  '
$ws.Range("D31").Value = '116.84'
$ws.Range("E31").Value = '  -2.64%  '
$ws.Range("D32").Value = '0.9511'
$ws.Range("E32").Value = '  -3.34%  '
$ws.Range("D33").Value = '0.09275'
$ws.Range("E33").Value = '  -2.24%  '
$ws.Range("D34").Value = '3.586'
$ws.Range("E34").Value = '  +0.52%  '
$ws.Range("D35").Value = '1.394'
$ws.Range("E35").Value = '  +0.99%  '
$ws.Range("D36").Value = '5.218'
$ws.Range("E36").Value = '  -2.87%  '
$ws.Range("D37").Value = '0.06004'
$ws.Range("E37").Value = '  -1.82%  '
$ws.Range("D38").Value = '0.02201'
$ws.Range("E38").Value = '  -2.95%  '
$ws.Range("D39").Value = '8.254'
$ws.Range("E39").Value = '  -2.39%  '
$ws.Range("D40").Value = '1.166'
$ws.Range("E40").Value = '  -1.15%  '
$ws.Range("D41").Value = '0.9957'
$ws.Range("E41").Value = '  -0.34%  '
$ws.Range("D42").Value = '0.5722'
$ws.Range("E42").Value = '  -4.37%  '
$ws.Range("B43").Value = 'Aptos'
$ws.Range("C43").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D43").Value = '10.03'
$ws.Range("E43").Value = '  -3.95%  '
$ws.Range("B44").Value = 'Algorand'
$ws.Range("C44").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D44").Value = '0.1800'
$ws.Range("E44").Value = '  -4.67%  '
$ws.Range("D45").Value = '1.236'
$ws.Range("E45").Value = '  -3.57%  '
$ws.Range("D46").Value = '2.318'
$ws.Range("E46").Value = '  +16.28%  '
$ws.Range("D47").Value = '0.5416'
$ws.Range("E47").Value = '  -4.03%  '
$ws.Range("D48").Value = '11.86'
$ws.Range("E48").Value = '  -3.31%  '
$ws.Range("D49").Value = '0.07197'
$ws.Range("E49").Value = '  +4.37%  '
$ws.Range("D50").Value = '1.869'
$ws.Range("E50").Value = '  -5.43%  '
$ws.Range("D51").Value = '109.64'
$ws.Range("E51").Value = '  -1.64%  '
